$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1064.4814
$ws.Cells.Item(17, 10).Value = 1064.4814
$ws.Cells.Item(17, 12).Value = 3193.4442
$ws.Cells.Item(17, 14).Value = -3529.4442

$ws.Cells.Item(54, 8).Value = 15000
$ws.Cells.Item(54, 9).Value = 15000
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 11).Value = 15000
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 13).Value = -14514
$ws.Cells.Item(54, 14).ClearContents()

$ws.Cells.Item(64, 8).Value = 8178.5356
$ws.Cells.Item(64, 10).Value = 9454.5
$ws.Cells.Item(64, 12).Value = 9454.5
$ws.Cells.Item(64, 14).Value = -9950.5

$ws.Cells.Item(67, 8).Value = 8178.5356
$ws.Cells.Item(67, 10).Value = 9454.5
$ws.Cells.Item(67, 12).Value = 9454.5
$ws.Cells.Item(67, 14).Value = -11170.5

$ws.Cells.Item(70, 8).Value = 13338436
$ws.Cells.Item(70, 10).Value = 9000.375
$ws.Cells.Item(70, 12).Value = 27001.125
$ws.Cells.Item(70, 14).Value = -27541.125

$ws.Cells.Item(73, 8).Value = 13338436
$ws.Cells.Item(73, 10).Value = 9000.375
$ws.Cells.Item(73, 12).Value = 27001.125
$ws.Cells.Item(73, 14).Value = -28873.125

$ws.Cells.Item(74, 8).Value = 13319.167
$ws.Cells.Item(74, 9).Value = 13319.167
$ws.Cells.Item(74, 11).Value = 13319.167
$ws.Cells.Item(74, 13).Value = -12383.167

$ws.Cells.Item(77, 8).Value = 13319.167
$ws.Cells.Item(77, 9).Value = 13319.167
$ws.Cells.Item(77, 11).Value = 66595.83499999999
$ws.Cells.Item(77, 13).Value = -61915.83499999999

$ws.Cells.Item(101, 8).Value = 427.25
$ws.Cells.Item(101, 9).Value = 427.25
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 1281.75
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 13).Value = 340.25
$ws.Cells.Item(101, 14).ClearContents()

$ws.Cells.Item(134, 8).Value = 50000
$ws.Cells.Item(134, 9).Value = 50000
$ws.Cells.Item(134, 11).Value = 50000
$ws.Cells.Item(134, 13).Value = -44930

$ws.Cells.Item(137, 8).Value = 1931.2941
$ws.Cells.Item(137, 9).Value = 1855.4667
$ws.Cells.Item(137, 11).Value = 5566.4001
$ws.Cells.Item(137, 13).Value = -3016.4001

$ws.Cells.Item(138, 8).Value = 1256.3715
$ws.Cells.Item(138, 9).Value = 1043.6875
$ws.Cells.Item(138, 11).Value = 3131.0625
$ws.Cells.Item(138, 13).Value = 2008.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()

$ws.Cells.Item(45, 8).Value = 1500
$ws.Cells.Item(45, 9).Value = 1111.1111
$ws.Cells.Item(45, 10).Value = 5000
$ws.Cells.Item(45, 11).Value = 1111.1111
$ws.Cells.Item(45, 12).Value = 5000
$ws.Cells.Item(45, 13).Value = -734.1111000000001
$ws.Cells.Item(45, 14).Value = -5754

$ws.Cells.Item(74, 8).Value = 1706.3334
$ws.Cells.Item(74, 9).Value = 1691.65
$ws.Cells.Item(74, 11).Value = 1691.65
$ws.Cells.Item(74, 13).Value = -817.6500000000001

$ws.Cells.Item(77, 8).Value = 1706.3334
$ws.Cells.Item(77, 9).Value = 1691.65
$ws.Cells.Item(77, 11).Value = 8458.25
$ws.Cells.Item(77, 13).Value = -4090.25

$ws.Cells.Item(132, 8).Value = 1641.2833
$ws.Cells.Item(132, 9).Value = 1672.9434
$ws.Cells.Item(132, 10).Value = 1401.5714
$ws.Cells.Item(132, 11).Value = 5018.8302
$ws.Cells.Item(132, 12).Value = 4204.7142
$ws.Cells.Item(132, 13).Value = -2488.8302
$ws.Cells.Item(132, 14).Value = -9264.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(62, 8).Value = 42500
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 42500
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 42500
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).Value = -43872

$ws.Cells.Item(65, 8).Value = 42500
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 42500
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 127500
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).Value = -134364

$ws.Cells.Item(99, 8).Value = 8000
$ws.Cells.Item(99, 10).Value = 8000
$ws.Cells.Item(99, 12).Value = 8000
$ws.Cells.Item(99, 14).Value = -10996

$ws.Cells.Item(107, 8).Value = 5100.231
$ws.Cells.Item(107, 9).Value = 4011
$ws.Cells.Item(107, 11).Value = 4011
$ws.Cells.Item(107, 13).Value = -2091

$ws.Cells.Item(126, 8).Value = 75000
$ws.Cells.Item(126, 10).Value = 75000
$ws.Cells.Item(126, 12).Value = 75000
$ws.Cells.Item(126, 14).Value = -84880

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 4334.75
$ws.Cells.Item(17, 9).Value = 3602.8
$ws.Cells.Item(17, 10).Value = 5554.6665
$ws.Cells.Item(17, 11).Value = 3602.8
$ws.Cells.Item(17, 12).Value = 5554.6665
$ws.Cells.Item(17, 13).Value = -3428.8
$ws.Cells.Item(17, 14).Value = -5902.6665

$ws.Cells.Item(22, 8).Value = 1371.6875
$ws.Cells.Item(22, 9).Value = 1230.2727
$ws.Cells.Item(22, 10).Value = 1682.8
$ws.Cells.Item(22, 11).Value = 1230.2727
$ws.Cells.Item(22, 12).Value = 1682.8
$ws.Cells.Item(22, 13).Value = -880.2727
$ws.Cells.Item(22, 14).Value = -2382.8

$ws.Cells.Item(58, 8).Value = 3498.923
$ws.Cells.Item(58, 9).Value = 2908
$ws.Cells.Item(58, 11).Value = 2908
$ws.Cells.Item(58, 13).Value = -2705

$ws.Cells.Item(122, 8).Value = 894.2353000000001
$ws.Cells.Item(122, 10).Value = 1089.25
$ws.Cells.Item(122, 12).Value = 3267.75
$ws.Cells.Item(122, 14).Value = -8167.75

$ws.Cells.Item(136, 8).Value = 3498.923
$ws.Cells.Item(136, 9).Value = 2908
$ws.Cells.Item(136, 11).Value = 8724
$ws.Cells.Item(136, 13).Value = -6174

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 7000
$ws.Cells.Item(20, 9).Value = 7000
$ws.Cells.Item(20, 11).Value = 21000
$ws.Cells.Item(20, 13).Value = -20773

$ws.Cells.Item(56, 8).Value = 48250.285
$ws.Cells.Item(56, 9).Value = 48250.285
$ws.Cells.Item(56, 11).Value = 48250.285
$ws.Cells.Item(56, 13).Value = -47720.285

$ws.Cells.Item(87, 8).Value = 3861
$ws.Cells.Item(87, 9).Value = 3861
$ws.Cells.Item(87, 11).Value = 11583
$ws.Cells.Item(87, 13).Value = -10335

$ws.Cells.Item(90, 8).Value = 3861
$ws.Cells.Item(90, 9).Value = 3861
$ws.Cells.Item(90, 11).Value = 34749
$ws.Cells.Item(90, 13).Value = -28509

$ws.Cells.Item(103, 8).Value = 1500
$ws.Cells.Item(103, 9).Value = 1500
$ws.Cells.Item(103, 11).Value = 4500
$ws.Cells.Item(103, 13).Value = -3621

$ws.Cells.Item(113, 8).Value = 1093.091
$ws.Cells.Item(113, 10).Value = 1799.8
$ws.Cells.Item(113, 12).Value = 5399.4
$ws.Cells.Item(113, 14).Value = -9739.4

$ws.Cells.Item(114, 8).Value = 857.5
$ws.Cells.Item(114, 10).Value = 1010
$ws.Cells.Item(114, 12).Value = 3030
$ws.Cells.Item(114, 14).Value = -9538

$ws.Cells.Item(125, 8).Value = 10000
$ws.Cells.Item(125, 9).Value = 10000
$ws.Cells.Item(125, 11).Value = 30000
$ws.Cells.Item(125, 13).Value = -25080

$ws.Cells.Item(126, 8).Value = 100
$ws.Cells.Item(126, 9).Value = 100
$ws.Cells.Item(126, 11).Value = 300
$ws.Cells.Item(126, 13).Value = 4640

$ws.Cells.Item(132, 8).Value = 2659.8667
$ws.Cells.Item(132, 9).Value = 2316.5
$ws.Cells.Item(132, 10).Value = 2888.7778
$ws.Cells.Item(132, 11).Value = 20848.5
$ws.Cells.Item(132, 12).Value = 25999.0002
$ws.Cells.Item(132, 13).Value = -18318.5
$ws.Cells.Item(132, 14).Value = -31059.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64, 8).Value = 75000
$ws.Cells.Item(64, 10).Value = 75000
$ws.Cells.Item(64, 12).Value = 75000
$ws.Cells.Item(64, 14).Value = -75496

$ws.Cells.Item(67, 8).Value = 75000
$ws.Cells.Item(67, 10).Value = 75000
$ws.Cells.Item(67, 12).Value = 75000
$ws.Cells.Item(67, 14).Value = -76716

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 2184.0833
$ws.Cells.Item(55, 9).Value = 343.2
$ws.Cells.Item(55, 10).Value = 3499
$ws.Cells.Item(55, 11).Value = 343.2
$ws.Cells.Item(55, 12).Value = 3499
$ws.Cells.Item(55, 13).Value = -170.2
$ws.Cells.Item(55, 14).Value = -3845
